$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price/volume text values stay as text (preserve exact formatting,
# e.g. trailing zeros) instead of being auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = "310.96"
$ws.Cells.Item(2, 5).Value = "-1.36%"
$ws.Cells.Item(3, 4).Value = "48.58"
$ws.Cells.Item(3, 5).Value = "8.67%"
$ws.Cells.Item(4, 4).Value = "5.250"
$ws.Cells.Item(4, 5).Value = "1.88%"
$ws.Cells.Item(5, 4).Value = "0.07772"
$ws.Cells.Item(5, 5).Value = "-3.56%"
$ws.Cells.Item(6, 4).Value = "4.532"
$ws.Cells.Item(6, 5).Value = "0.31%"
$ws.Cells.Item(7, 4).Value = "1.301"
$ws.Cells.Item(7, 5).Value = "19.62%"
$ws.Cells.Item(8, 4).Value = "1.561"
$ws.Cells.Item(8, 5).Value = "-7.10%"
$ws.Cells.Item(9, 4).Value = "0.1246"
$ws.Cells.Item(9, 5).Value = "-3.77%"
$ws.Cells.Item(10, 4).Value = "0.1951"
$ws.Cells.Item(10, 5).Value = "1.75%"
$ws.Cells.Item(11, 4).Value = "0.09314"
$ws.Cells.Item(11, 5).Value = "-0.52%"
$ws.Cells.Item(12, 4).Value = "0.04565"
$ws.Cells.Item(12, 5).Value = "7.42%"
$ws.Cells.Item(14, 4).Value = "0.001300"
$ws.Cells.Item(14, 5).Value = "-1.57%"
$ws.Cells.Item(15, 4).Value = "0.04212"
$ws.Cells.Item(15, 5).Value = "0.00%"
$ws.Cells.Item(16, 4).Value = "0.005863"
$ws.Cells.Item(16, 5).Value = "-1.12%"
$ws.Cells.Item(17, 2).Value = "LEO"
$ws.Cells.Item(17, 3).Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Cells.Item(17, 4).Value = "3.327"
$ws.Cells.Item(17, 5).Value = "-1.89%"
$ws.Cells.Item(18, 2).Value = "BTSEToken"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Cells.Item(18, 4).Value = "2.410"
$ws.Cells.Item(18, 5).Value = "0.04%"
$ws.Cells.Item(19, 2).Value = "BitpandaEcosystemToken"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Cells.Item(19, 4).Value = "0.3443"
$ws.Cells.Item(19, 5).Value = "1.51%"
$ws.Cells.Item(20, 2).Value = "MCDex"
$ws.Cells.Item(20, 3).Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Cells.Item(20, 4).Value = "8.135"
$ws.Cells.Item(20, 5).Value = "-1.68%"
$ws.Cells.Item(21, 2).Value = "ProBitToken"
$ws.Cells.Item(21, 3).Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Cells.Item(21, 4).Value = "0.1372"
$ws.Cells.Item(21, 5).Value = "-0.71%"
$ws.Cells.Item(22, 2).Value = "ZBToken"
$ws.Cells.Item(22, 3).Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Cells.Item(22, 4).Value = "0.3069"
$ws.Cells.Item(22, 5).Value = "-2.18%"
$ws.Cells.Item(23, 2).Value = "BitKan"
$ws.Cells.Item(23, 3).Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Cells.Item(23, 4).Value = "0.001299"
$ws.Cells.Item(23, 5).Value = "1.96%"
$ws.Cells.Item(24, 2).Value = "HotbitToken"
$ws.Cells.Item(24, 3).Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Cells.Item(24, 4).Value = "0.004153"
$ws.Cells.Item(24, 5).Value = "-8.78%"
$ws.Cells.Item(25, 4).Value = "0.0001360"
$ws.Cells.Item(25, 5).Value = "1.38%"
$ws.Cells.Item(26, 4).Value = "0.0003567"
$ws.Cells.Item(26, 5).Value = "-95.19%"
$ws.Cells.Item(38, 4).Value = "0.02571"
$ws.Cells.Item(38, 5).Value = "-4.64%"
$ws.Cells.Item(39, 4).Value = "0.05804"
$ws.Cells.Item(39, 5).Value = "6.65%"
$ws.Cells.Item(40, 4).Value = "0.01041"
$ws.Cells.Item(40, 5).Value = "87.07%"
$ws.Cells.Item(41, 4).Value = "0.007978"
$ws.Cells.Item(41, 5).Value = "3.17%"
$ws.Cells.Item(42, 4).Value = "0.1418"
$ws.Cells.Item(42, 5).Value = "-0.14%"
$ws.Cells.Item(43, 4).Value = "0.008448"
$ws.Cells.Item(43, 5).Value = "15.36%"
$ws.Cells.Item(44, 4).Value = "0.008525"
$ws.Cells.Item(44, 5).Value = "-0.53%"
$ws.Cells.Item(45, 4).Value = "0.3121"
$ws.Cells.Item(45, 5).Value = "-0.56%"
$ws.Cells.Item(46, 4).Value = "0.00006931"
$ws.Cells.Item(46, 5).Value = "2.07%"
$ws.Cells.Item(47, 4).Value = "0.00000000756"
$ws.Cells.Item(47, 5).Value = "1.39%"
$ws.Cells.Item(48, 4).Value = "0.05500"
$ws.Cells.Item(48, 5).Value = "-10.65%"
$ws.Cells.Item(49, 4).Value = "0.004031"
$ws.Cells.Item(49, 5).Value = "1.39%"
$ws.Cells.Item(50, 4).Value = "0.00002116"
$ws.Cells.Item(50, 5).Value = "1.39%"
$ws.Cells.Item(51, 4).Value = "0.0002016"
$ws.Cells.Item(51, 5).Value = "1.39%"
